$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally holds two lists side by side in columns D (names) and E (emails),
# rows 1-8 hold both, rows 9-28 only hold emails in column E.
# The edit inserts three new plain numeric columns (D, E, F) before the existing lists,
# pushing the "TONAME"/"EMAIL" lists out to columns G and H.

# Move column D (which has data only in rows 1-8) to column G.
$ws.Range("D1:D8").Copy($ws.Range("G1:G8"))

# Move column E (which has data in rows 1-28) to column H.
$ws.Range("E1:E28").Copy($ws.Range("H1:H28"))

# Clear the old D:E content/formatting now that it has been relocated.
$ws.Range("D1:E28").Clear()

# Populate the three new columns in the header row with plain numbers.
$ws.Range("D1").Value = 4
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 6

# Match the resulting selection shown in the saved workbook.
$ws.Range("F2").Select()
